# Generate Report for Handback
#
# For the "97cdfdf1-cd49-46ea-b470-466da18b27ac" row (row 8) on both the
# zh-cn and de-de sheets, the handback processing has now produced a
# target/handback file, a handback datetime and an error (the handed back
# file version is not the latest available), so the previously-empty
# "Latest Target File", "Latest Handback File", "Latest Handback DateTime"
# and "Error Detail" columns get populated. The "Error Detail" column is
# also widened so the long message is readable.

$wb = $excel.ActiveWorkbook

$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bb62e1dfadc09451ba7548c290f07354364f4fe/e2e/97cdfdf1-cd49-46ea-b470-466da18b27ac.md"
$latestUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a2b92bf4d878826acb1544d2e33e093ea2858bc/e2e/97cdfdf1-cd49-46ea-b470-466da18b27ac.md"
$targetMd   = "97cdfdf1-cd49-46ea-b470-466da18b27ac.md"
$errorMsg   = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."

# Colour used by this workbook's existing "HyperLink" style (font color FF6495ED, BGR-encoded for Font.Color).
$hyperlinkColor = 15570276

function Update-LocalizationSheet($Sheet, $HandbackXlf, $HandbackDateTime) {
    # Widen the "Error Detail" column (P) so the new message is legible.
    $Sheet.Range("P1").EntireColumn.ColumnWidth = 39.1667

    # Latest Target File (I8): becomes a hyperlink to the (non-latest) handed
    # back version of the source markdown file.
    $Sheet.Hyperlinks.Add($Sheet.Range("I8"), $currentUrl, "", "", $targetMd)
    $Sheet.Range("I8").Font.Color = $hyperlinkColor
    $Sheet.Range("I8").Font.Underline = 2

    # Latest Handback File (J8)
    $Sheet.Range("J8").Value = $HandbackXlf

    # Latest Handback DateTime (K8)
    $Sheet.Range("K8").Value = $HandbackDateTime

    # Error Detail (P8)
    $Sheet.Range("P8").Value = $errorMsg
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Update-LocalizationSheet $zhcn "97cdfdf1-cd49-46ea-b470-466da18b27ac.fbd8c4df2cef5103de5583ce8d6fd432d7850341.zh-cn.xlf" "2016-08-26 08:47:28"

$dede = $wb.Worksheets.Item("de-de")
Update-LocalizationSheet $dede "97cdfdf1-cd49-46ea-b470-466da18b27ac.fbd8c4df2cef5103de5583ce8d6fd432d7850341.de-de.xlf" "2016-08-26 08:47:35"
